# Actualización automática 2025-07-28 15:20:09
# Inserts a new advisor "JIMENEZ CORDERO WILLIAM GUSTAVO" (with all-zero
# figures) ahead of "LOZANO MOLINA TITO JERSON" in both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, shifting the existing
# rows down by one and updating the trailing "N de 14" counters to
# "N de 15" on the "VENTAS POR GRUPO" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: VENTAS POR GRUPO  (columns A:R, new row inserted at row 9)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(9).Insert()

$ws1.Cells.Item(9, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(9, 2).Value = "JIMENEZ CORDERO WILLIAM GUSTAVO"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(9, $c).Value = 0
}

# the summary row (now row 17) shows counters like "0 de 14" — bump the
# denominator now that there are 15 advisors listed instead of 14
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(17, $c)
    $cell.Value = $cell.Value2().Replace("de 14", "de 15")
}

# ---------------------------------------------------------------
# Sheet: VENTA MENSUAL  (columns A:G, new row inserted at row 9)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(9).Insert()

$ws2.Cells.Item(9, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(9, 2).Value = "JIMENEZ CORDERO WILLIAM GUSTAVO"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(9, $c).Value = 0
}
